$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2371.7163
$ws.Range("I15").Value = 2371.7163
$ws.Range("K15").Value = 7115.1489
$ws.Range("M15").Value = -6946.1489
$ws.Range("H28").Value = 7147.8335
$ws.Range("I28").Value = 465.5
$ws.Range("J28").Value = 15500.75
$ws.Range("K28").Value = 465.5
$ws.Range("L28").Value = 15500.75
$ws.Range("M28").Value = 19.5
$ws.Range("N28").Value = -16470.75
$ws.Range("H107").Value = 966.2
$ws.Range("I107").Value = 902.4091
$ws.Range("J107").Value = 1434
$ws.Range("K107").Value = 902.4091
$ws.Range("L107").Value = 1434
$ws.Range("M107").Value = 1017.5909
$ws.Range("N107").Value = -5274
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 685364.4
$ws.Range("I32").Value = 896337.7
$ws.Range("J32").Value = 14085.637
$ws.Range("K32").Value = 896337.7
$ws.Range("L32").Value = 14085.637
$ws.Range("M32").Value = -896050.7
$ws.Range("N32").Value = -14659.637
$ws.Range("H45").Value = 3478
$ws.Range("I45").Value = 2473.1428
$ws.Range("J45").Value = 4181.4
$ws.Range("K45").Value = 2473.1428
$ws.Range("L45").Value = 4181.4
$ws.Range("M45").Value = -2096.1428
$ws.Range("N45").Value = -4935.4
$ws.Range("H74").Value = 1717.7142
$ws.Range("I74").Value = 1542.1538
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 1542.1538
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -668.1538
$ws.Range("N74").Value = -5748
$ws.Range("H77").Value = 1717.7142
$ws.Range("I77").Value = 1542.1538
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 7710.769
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -3342.769
$ws.Range("N77").Value = -28736
$ws.Range("H101").Value = 79602
$ws.Range("J101").Value = 79602
$ws.Range("L101").Value = 79602
$ws.Range("N101").Value = -86092
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 11780
$ws.Range("J57").Value = 11780
$ws.Range("L57").Value = 11780
$ws.Range("N57").Value = -13220
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H136").Value = 11780
$ws.Range("J136").Value = 11780
$ws.Range("L136").Value = 11780
$ws.Range("N136").Value = -21980
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 975
$ws.Range("I16").Value = 900
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -613
$ws.Range("N16").Value = -1574
$ws.Range("H113").Value = 975
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -5340
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 113
$ws.Range("J2").Value = 196
$ws.Range("L2").Value = 1176
$ws.Range("N2").Value = -1402
$ws.Range("H5").Value = 1024.5
$ws.Range("I5").Value = 899.4
$ws.Range("K5").Value = 2698.2
$ws.Range("M5").Value = -2586.2
$ws.Range("H26").Value = 6398.697
$ws.Range("I26").Value = 32.5
$ws.Range("J26").Value = 8435.879999999999
$ws.Range("K26").Value = 97.5
$ws.Range("L26").Value = 25307.64
$ws.Range("M26").Value = 190.5
$ws.Range("N26").Value = -25883.64
$ws.Range("H34").Value = 9804338
$ws.Range("J34").Value = 10417103
$ws.Range("L34").Value = 31251309
$ws.Range("N34").Value = -31251477
$ws.Range("H39").Value = 1690.5834
$ws.Range("I39").Value = 290
$ws.Range("J39").Value = 1817.909
$ws.Range("K39").Value = 870
$ws.Range("L39").Value = 5453.727000000001
$ws.Range("M39").Value = -576
$ws.Range("N39").Value = -6041.727000000001
$ws.Range("H75").Value = 5307.5
$ws.Range("I75").Value = 3080
$ws.Range("J75").Value = 6898.5713
$ws.Range("K75").Value = 9240
$ws.Range("L75").Value = 20695.7139
$ws.Range("M75").Value = -8242
$ws.Range("N75").Value = -22691.7139
$ws.Range("H78").Value = 5307.5
$ws.Range("I78").Value = 3080
$ws.Range("J78").Value = 6898.5713
$ws.Range("K78").Value = 27720
$ws.Range("L78").Value = 62087.14169999999
$ws.Range("M78").Value = -22728
$ws.Range("N78").Value = -72071.14169999999
$ws.Range("H113").Value = 1542.75
$ws.Range("I113").Value = 850
$ws.Range("J113").Value = 1889.125
$ws.Range("K113").Value = 2550
$ws.Range("L113").Value = 5667.375
$ws.Range("M113").Value = -380
$ws.Range("N113").Value = -10007.375
$ws.Range("H122").Value = 7722.9287
$ws.Range("I122").Value = 440.5
$ws.Range("J122").Value = 17432.834
$ws.Range("K122").Value = 3964.5
$ws.Range("L122").Value = 156895.506
$ws.Range("M122").Value = -1514.5
$ws.Range("N122").Value = -161795.506
$ws.Range("H135").Value = 1024.5
$ws.Range("I135").Value = 899.4
$ws.Range("K135").Value = 8094.599999999999
$ws.Range("M135").Value = -5559.599999999999
$ws.Range("H137").Value = 20848176
$ws.Range("I137").Value = 20848176
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 62544528
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -62539428
$ws.Range("N137").ClearContents()
$ws.Range("H140").Value = 2057.0715
$ws.Range("I140").Value = 1481.7273
$ws.Range("J140").Value = 4166.6665
$ws.Range("K140").Value = 4445.1819
$ws.Range("L140").Value = 12499.9995
$ws.Range("M140").Value = 734.8181000000004
$ws.Range("N140").Value = -22859.9995
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 29999.334
$ws.Range("J15").Value = 29999.334
$ws.Range("L15").Value = 29999.334
$ws.Range("N15").Value = -30575.334
$ws.Range("H80").Value = 450765.56
$ws.Range("I80").Value = 603204.9399999999
$ws.Range("J80").Value = 69667.164
$ws.Range("K80").Value = 603204.9399999999
$ws.Range("L80").Value = 69667.164
$ws.Range("M80").Value = -602206.9399999999
$ws.Range("N80").Value = -71663.164
$ws.Range("H81").Value = 29999.334
$ws.Range("J81").Value = 29999.334
$ws.Range("L81").Value = 29999.334
$ws.Range("N81").Value = -31995.334
$ws.Range("H83").Value = 450765.56
$ws.Range("I83").Value = 603204.9399999999
$ws.Range("J83").Value = 69667.164
$ws.Range("K83").Value = 3016024.7
$ws.Range("L83").Value = 348335.82
$ws.Range("M83").Value = -3011032.7
$ws.Range("N83").Value = -358319.82
$ws.Range("H84").Value = 29999.334
$ws.Range("J84").Value = 29999.334
$ws.Range("L84").Value = 89998.00199999999
$ws.Range("N84").Value = -99982.00199999999
$ws.Range("H101").Value = 49379.57
$ws.Range("J101").Value = 49379.57
$ws.Range("L101").Value = 49379.57
$ws.Range("N101").Value = -55869.57
$ws.Range("H107").Value = 707.0769
$ws.Range("I107").Value = 648.4545000000001
$ws.Range("J107").Value = 750.06665
$ws.Range("K107").Value = 648.4545000000001
$ws.Range("L107").Value = 750.06665
$ws.Range("M107").Value = 1271.5455
$ws.Range("N107").Value = -4590.06665
$ws.Range("H132").Value = 2175.3103
$ws.Range("I132").Value = 1694.3158
$ws.Range("K132").Value = 5082.9474
$ws.Range("M132").Value = -2552.9474
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 52712.832
$ws.Range("J140").Value = 52712.832
$ws.Range("L140").Value = 52712.832
$ws.Range("N140").Value = -63072.832
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 18334.5
$ws.Range("J15").Value = 18334.5
$ws.Range("L15").Value = 18334.5
$ws.Range("N15").Value = -18910.5
$ws.Range("H31").Value = 68389.25
$ws.Range("J31").Value = 68389.25
$ws.Range("L31").Value = 68389.25
$ws.Range("N31").Value = -69085.25
